# Hortaliza, Vega Monumental Concepción - Brócoli
# A new weekly data point is inserted as row 88 (Brócoli / Primera / Vega
# Monumental Concepción), pushing the existing rows 88-187 down to 89-188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 88, shifting rows 88:187 down to 89:188.
$ws.Rows("88").Insert()

# Populate the newly inserted row 88 with the new observation.
$ws.Range("A88").Value = 11
$ws.Range("B88").Value = "Vega Monumental Concepción"
$ws.Range("C88").Value = "Bíobío"
$ws.Range("D88").Value = 44539
$ws.Range("E88").Value = 8
$ws.Range("F88").Value = 100112023
$ws.Range("G88").Value = "Brócoli"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 1900
$ws.Range("K88").Value = 500
$ws.Range("L88").Value = 600
$ws.Range("M88").Value = 547
$ws.Range("N88").Value = "$/unidad"
$ws.Range("O88").Value = "Región Metropolitana"
$ws.Range("P88").Value = 547
$ws.Range("Q88").Value = 1
$ws.Range("R88").Value = "Hortaliza"
